$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.208.38"
$ws.Range("E2").Value = "  -2.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.820.82"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "422.92"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.31"
$ws.Range("E6").Value = "  -3.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.818.24"
$ws.Range("E7").Value = "  +1.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("E8").Value = "  -7.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.716"
$ws.Range("E10").Value = "  -7.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  -13.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000345"
$ws.Range("E12").Value = "  -20.86%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.01"
$ws.Range("E13").Value = "  -7.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.409.87"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.89"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.69"
$ws.Range("E16").Value = "  +19.67%  "

$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.824.09"
$ws.Range("E18").Value = "  +1.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.44"
$ws.Range("E19").Value = "  -6.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "66.387.23"
$ws.Range("E20").Value = "  -2.33%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  -7.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "401.01"
$ws.Range("E22").Value = "  -11.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.24"
$ws.Range("E23").Value = "  -10.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.56"
$ws.Range("E24").Value = "  -6.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.97"
$ws.Range("E25").Value = "  -4.17%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "36.89"
$ws.Range("E26").Value = "  -5.17%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.76"
$ws.Range("E27").Value = "  +12.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.18"
$ws.Range("E28").Value = "  -4.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.38"
$ws.Range("E29").Value = "  -7.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "693.87"
$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.120"
$ws.Range("E31").Value = "  -4.98%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.21"
$ws.Range("E33").Value = "  -4.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.33"
$ws.Range("E34").Value = "  +1.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.149"
$ws.Range("E35").Value = "  -10.29%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.67"
$ws.Range("E36").Value = "  -10.37%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "54.67"
$ws.Range("E38").Value = "  -4.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0758"
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0449"
$ws.Range("E40").Value = "  -9.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -10.49%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.43"
$ws.Range("E44").Value = "  +1.36%  "

$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -4.25%  "

$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.85"
$ws.Range("E46").Value = "  -2.64%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.06"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.99"
$ws.Range("E48").Value = "  -6.40%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("E49").Value = "  -5.12%  "

$ws.Range("E50").Value = "  -4.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.71"
$ws.Range("E51").Value = "  -7.40%  "
